# worked on auto grapher, still trash, dont use
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear K4 (value 56 removed from sheet data)
$ws.Range("K4").ClearContents()

# Update the active selection to Q5
$ws.Range("Q5").Select()
